# Cinch V1 Release Notes - apply release note edits
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Heading: "Release XXXX" -> "Current Release : 03/05/2010 10:30PM Uk time"
#    (rebuilt as multiple runs, with "Uk" flagged as a spelling error, to
#    mirror the original author's run/proofErr layout)
# ---------------------------------------------------------------------------
$headingXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/></w:rPr><w:t xml:space="preserve">Current </w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/></w:rPr><w:t>Release</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/></w:rPr><w:t xml:space="preserve"> : 03/05/2010 10:30PM </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/></w:rPr><w:t>Uk</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/></w:rPr><w:t xml:space="preserve"> time</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r = $d.Content
$r.Find.Execute("Release XXXX", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r.InsertXML($headingXml)

# ---------------------------------------------------------------------------
# 2) Table work: shade + vertically merge the "Notes" column for the three
#    Mediator-related rows, and rewrite the cell text content.
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)

# -- Row 2 / Col 1 ("Put in new Mediator, and make ViewModel...") --
$cell21 = $t.Cell(2, 1)
$cell21.Shading.Texture = 0
$cell21.Shading.ForegroundPatternColor = -16777216
$cell21.Shading.BackgroundPatternColor = 5296274
$cell21Xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve">Put in new Mediator, and make </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>ViewModel</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> unregister on Dispose</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$cell21.Range.InsertXML($cell21Xml)

# -- Row 2 / Col 2 (new explanatory notes) --
$cell22 = $t.Cell(2, 2)
$cell22Xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve">There has been a lot of work done in the Mediator which is now a Singleton, which can be used within Views also, and also extra Unregister/Register methods have been made available, as well as </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>Async</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> methods.</w:t></w:r>
</w:p>
<w:p/>
<w:p>
<w:r><w:t xml:space="preserve">The </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:t>ViewModelBase</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:t xml:space="preserve"> also unregisters within the </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>Dispose(</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t>) method.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$cell22.Range.InsertXML($cell22Xml)

# -- Shade col-1 of rows 3 and 4 (no content changes there) --
$cell31 = $t.Cell(3, 1)
$cell31.Shading.Texture = 0
$cell31.Shading.ForegroundPatternColor = -16777216
$cell31.Shading.BackgroundPatternColor = 5296274

$cell41 = $t.Cell(4, 1)
$cell41.Shading.Texture = 0
$cell41.Shading.ForegroundPatternColor = -16777216
$cell41.Shading.BackgroundPatternColor = 5296274

# -- Shade col-2 of rows 2, 3 and 4, then vertically merge them --
# (shade while each row is still individually addressable; after Merge()
#  only the top cell of the merged range stays reachable via Cell(row,col))
$t2 = $d.Tables.Item(1)
$cell32 = $t2.Cell(3, 2)
$cell32.Shading.Texture = 0
$cell32.Shading.ForegroundPatternColor = -16777216
$cell32.Shading.BackgroundPatternColor = 5296274

$cell42 = $t2.Cell(4, 2)
$cell42.Shading.Texture = 0
$cell42.Shading.ForegroundPatternColor = -16777216
$cell42.Shading.BackgroundPatternColor = 5296274

$cell22b = $t2.Cell(2, 2)
$cell22b.Shading.Texture = 0
$cell22b.Shading.ForegroundPatternColor = -16777216
$cell22b.Shading.BackgroundPatternColor = 5296274

$t3 = $d.Tables.Item(1)
$mergeTop = $t3.Cell(2, 2)
$mergeBottom = $t3.Cell(4, 2)
$mergeTop.Merge($mergeBottom)

# ---------------------------------------------------------------------------
# 3) Add "46009 and " right before "Older Releases"
# ---------------------------------------------------------------------------
$olderXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/></w:rPr><w:t xml:space="preserve">46009 and </w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$olderRange = $d.Content
$olderRange.Find.Execute("Older ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertPoint = $d.Range($olderRange.Start, $olderRange.Start)
$insertPoint.InsertXML($olderXml)

Write-Host "Edits applied"
